$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) for rows 2-10 from serial date 45204 to 45207
foreach ($row in 2..10) {
    $ws.Cells.Item($row, 3).Value = 45207
}
